$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Date" column (column C) entirely, shifting remaining
# columns left.
$ws.Range("C1").EntireColumn.Delete()

# Remove the "OT/IT" column (originally column I, now column H after the
# previous delete) entirely, shifting remaining columns left.
$ws.Range("H1").EntireColumn.Delete()

# Update the sheet view's selection/scroll position to match the new
# layout.
$ws.Range("G26").Select()
